$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '255.93'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-0.78%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '26.99'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '0.32%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '4.363'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-8.14%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05881'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-1.59%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.631'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-0.83%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8516'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-2.55%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9328'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-2.34%'
$ws.Range("B9").Value = 'WazirX'
$ws.Range("C9").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1381'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-2.38%'
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.04699'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '29.86%'
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07083'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-2.22%'
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03069'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-2.13%'
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09096'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-1.57%'
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001527'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-1.72%'
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.006198'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '2.56%'
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.484'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-0.11%'
$ws.Range("B17").Value = 'GateToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.169'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-1.15%'
$ws.Range("B18").Value = 'BTSEToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.204'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-0.65%'
$ws.Range("B19").Value = 'One'
$ws.Range("C19").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.01036'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '1,606.52%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1270'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-1.60%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.916'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '10.71%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04271'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '1.31%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001218'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '0.37%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004280'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-5.08%'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-0.04%'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0001524'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '2.00%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03810'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-0.55%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006231'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '0.77%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1098'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-0.43%'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-2.34%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.01393'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '32.33%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005368'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-2.34%'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-0.04%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05501'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-35.70%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.2523'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '11,737.83%'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-0.04%'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.04%'
